$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 3999

$ws.Range("B4").Value = 3000
$ws.Range("C4").Value = 477
$ws.Range("D4").Value = 633
